$wb = $excel.ActiveWorkbook

# The sheet "會計科目表" holds the main account mapping table.
# Row 6 (NetIncomeLoss row) has the India-company columns (G:L, and N for
# Tech Mahindra) referencing the XBRL tag "Income". Correct it to the
# proper tag "ProfitLossForPeriod" to fix the India company ratio issue.
$ws1 = $wb.Worksheets.Item("會計科目表")

$ws1.Range("G6").Value = "ProfitLossForPeriod"
$ws1.Range("H6").Value = "ProfitLossForPeriod"
$ws1.Range("I6").Value = "ProfitLossForPeriod"
$ws1.Range("J6").Value = "ProfitLossForPeriod"
$ws1.Range("K6").Value = "ProfitLossForPeriod"
$ws1.Range("L6").Value = "ProfitLossForPeriod"
$ws1.Range("N6").Value = "ProfitLossForPeriod"

# Make "會計科目表" the active sheet (it was "會計科目表說明" before).
$ws1.Activate()
$ws1.Range("I12").Select()
